# This workbook's data table (rows 2-16, columns A:R) got re-shuffled:
# each destination row ends up holding the values that used to live in a
# different source row (a handful of rows - 4, 13, 14 - stay in place).
# Capture the original rows first, then write them back out in the new
# order so that no values are lost while we overwrite the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 16
$firstCol = 1   # A
$lastCol = 18   # R

# Mapping: destination row -> source row (values that should end up there)
$rowMap = @{
    2  = 8
    3  = 15
    4  = 4
    5  = 3
    6  = 12
    7  = 11
    8  = 16
    9  = 6
    10 = 2
    11 = 9
    12 = 10
    13 = 13
    14 = 14
    15 = 5
    16 = 7
}

# Snapshot all the original row values before writing anything back.
# (.Value2 is used because it reliably returns/accepts a plain 2-D array of
# raw values - dates as serials, numbers as numbers, strings as strings -
# for multi-cell ranges in this environment.)
$originalRows = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $originalRows[$r] = $rng.Value2
}

# Now write each destination row using the snapshot of its source row.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $srcRow = $rowMap[$r]
    $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $rng.Value2 = $originalRows[$srcRow]
}
